$d = $word.ActiveDocument

# Update the date line at the top of the document. (wdReplaceOne = 1, since
# ReplaceAll would re-scan/replace across the whole document regardless of
# the supplied range -- see note below.)
$d.Content.Find.Execute("2025-02-22 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-23 Sunday", 1)

# Update the answer table. Each non-blank table row holds 5 cells with
# "a÷b=c, d" style answers. Replacements are scoped to individual cells
# (re-fetched from the table each time) so duplicate text values in
# different cells are not confused with one another.
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "97÷3=32, 1";  New = "38÷8=4, 6" },
    @{ Row = 1;  Col = 2; Old = "18÷3=6, 0";   New = "12÷3=4, 0" },
    @{ Row = 1;  Col = 3; Old = "45÷9=5, 0";   New = "80÷3=26, 2" },
    @{ Row = 1;  Col = 4; Old = "53÷6=8, 5";   New = "90÷5=18, 0" },
    @{ Row = 1;  Col = 5; Old = "43÷8=5, 3";   New = "52÷6=8, 4" },

    @{ Row = 5;  Col = 1; Old = "58÷3=19, 1";  New = "39÷7=5, 4" },
    @{ Row = 5;  Col = 2; Old = "14÷6=2, 2";   New = "24÷8=3, 0" },
    @{ Row = 5;  Col = 3; Old = "45÷4=11, 1";  New = "48÷4=12, 0" },
    @{ Row = 5;  Col = 4; Old = "76÷4=19, 0";  New = "76÷2=38, 0" },
    @{ Row = 5;  Col = 5; Old = "77÷9=8, 5";   New = "73÷8=9, 1" },

    @{ Row = 9;  Col = 1; Old = "38÷6=6, 2";   New = "10÷5=2, 0" },
    @{ Row = 9;  Col = 2; Old = "64÷3=21, 1";  New = "28÷5=5, 3" },
    @{ Row = 9;  Col = 3; Old = "54÷7=7, 5";   New = "45÷5=9, 0" },
    @{ Row = 9;  Col = 4; Old = "34÷5=6, 4";   New = "45÷7=6, 3" },
    @{ Row = 9;  Col = 5; Old = "34÷4=8, 2";   New = "88÷5=17, 3" },

    @{ Row = 13; Col = 1; Old = "53÷3=17, 2";  New = "61÷7=8, 5" },
    @{ Row = 13; Col = 2; Old = "21÷7=3, 0";   New = "92÷3=30, 2" },
    @{ Row = 13; Col = 3; Old = "72÷4=18, 0";  New = "81÷3=27, 0" },
    @{ Row = 13; Col = 4; Old = "57÷8=7, 1";   New = "51÷3=17, 0" },
    @{ Row = 13; Col = 5; Old = "16÷9=1, 7";   New = "48÷6=8, 0" },

    @{ Row = 17; Col = 1; Old = "79÷9=8, 7";   New = "80÷2=40, 0" },
    @{ Row = 17; Col = 2; Old = "94÷5=18, 4";  New = "35÷8=4, 3" },
    @{ Row = 17; Col = 3; Old = "34÷4=8, 2";   New = "89÷4=22, 1" },
    @{ Row = 17; Col = 4; Old = "31÷3=10, 1";  New = "18÷2=9, 0" },
    @{ Row = 17; Col = 5; Old = "66÷3=22, 0";  New = "91÷2=45, 1" }
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    # Use wdReplaceOne (1) rather than wdReplaceAll (2): ReplaceAll is applied
    # document-wide (ignoring the Range boundary), which would wrongly touch
    # other cells that happen to share the same old text.
    $cell.Range.Find.Execute($c.Old, $true, $false, $false, $false, $false, $true, 1, $false, $c.New, 1)
}
